$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate column A (code) for new rows - creates shared string "html" (index 23)
$ws.Range("A11").Value = "html"
$ws.Range("A12").Value = "html"
$ws.Range("A13").Value = "html"

# Populate column B (descr) in the order that reproduces the target shared-string table:
# html(23), ملف html(24), html file(25), Fichier html(26)
$ws.Range("B12").Value = "ملف html"
$ws.Range("B11").Value = "html file"
$ws.Range("B13").Value = "Fichier html"

# Populate column C (lang_code) - reuse existing shared strings eng/ara/fra
$ws.Range("C11").Value = "eng"
$ws.Range("C12").Value = "ara"
$ws.Range("C13").Value = "fra"

# Populate column D (is_active) as boolean TRUE, with same style as other rows
$ws.Range("D11").Value = $true
$ws.Range("D12").Value = $true
$ws.Range("D13").Value = $true
$ws.Range("D11:D13").HorizontalAlignment = -4131

# Populate column E (cr_by) - reuse "superadmin"
$ws.Range("E11").Value = "superadmin"
$ws.Range("E12").Value = "superadmin"
$ws.Range("E13").Value = "superadmin"

# Populate column F (cr_dtimes) - reuse "now()"
$ws.Range("F11").Value = "now()"
$ws.Range("F12").Value = "now()"
$ws.Range("F13").Value = "now()"

# Update selection to match final state (full-column selection starting at G1)
$ws.Range("G1:XFD1048576").Select() | Out-Null
